# Apply the "time_taken" metadata column to the Heterotaxy panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamps for F2:F57, taken verbatim from the target diff.
$timestamps = @(
    "2021-10-05 10:51:04.382771",
    "2021-10-05 10:51:04.382782",
    "2021-10-05 10:51:04.382786",
    "2021-10-05 10:51:04.382789",
    "2021-10-05 10:51:04.382792",
    "2021-10-05 10:51:04.382795",
    "2021-10-05 10:51:04.382797",
    "2021-10-05 10:51:04.382800",
    "2021-10-05 10:51:04.382803",
    "2021-10-05 10:51:04.382806",
    "2021-10-05 10:51:04.382809",
    "2021-10-05 10:51:04.382812",
    "2021-10-05 10:51:04.382814",
    "2021-10-05 10:51:04.382817",
    "2021-10-05 10:51:04.382820",
    "2021-10-05 10:51:04.382822",
    "2021-10-05 10:51:04.382826",
    "2021-10-05 10:51:04.382828",
    "2021-10-05 10:51:04.382831",
    "2021-10-05 10:51:04.382834",
    "2021-10-05 10:51:04.382836",
    "2021-10-05 10:51:04.382839",
    "2021-10-05 10:51:04.382841",
    "2021-10-05 10:51:04.382844",
    "2021-10-05 10:51:04.382847",
    "2021-10-05 10:51:04.382850",
    "2021-10-05 10:51:04.382853",
    "2021-10-05 10:51:04.382855",
    "2021-10-05 10:51:04.382858",
    "2021-10-05 10:51:04.382860",
    "2021-10-05 10:51:04.382863",
    "2021-10-05 10:51:04.382866",
    "2021-10-05 10:51:04.382869",
    "2021-10-05 10:51:04.382871",
    "2021-10-05 10:51:04.382874",
    "2021-10-05 10:51:04.382876",
    "2021-10-05 10:51:04.382879",
    "2021-10-05 10:51:04.382882",
    "2021-10-05 10:51:04.382884",
    "2021-10-05 10:51:04.382887",
    "2021-10-05 10:51:04.382890",
    "2021-10-05 10:51:04.382892",
    "2021-10-05 10:51:04.382895",
    "2021-10-05 10:51:04.382898",
    "2021-10-05 10:51:04.382901",
    "2021-10-05 10:51:04.382903",
    "2021-10-05 10:51:04.382906",
    "2021-10-05 10:51:04.382908",
    "2021-10-05 10:51:04.382911",
    "2021-10-05 10:51:04.382913",
    "2021-10-05 10:51:04.382916",
    "2021-10-05 10:51:04.382918",
    "2021-10-05 10:51:04.382922",
    "2021-10-05 10:51:04.382924",
    "2021-10-05 10:51:04.382927",
    "2021-10-05 10:51:04.382929"
)

# Header cell F1, styled like the other header cells (B1:E1) by copying
# the formatting from the adjacent "panel" header cell (E1).
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data cells F2:F57 (no special style, matching the other data columns).
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
